# LOB1267.docx reshuffle:
# The edit rearranges which chunk of text occupies each paragraph slot
# (paragraph formatting / run formatting at every position is unchanged).
# Several of the new values are old values that currently live in other
# paragraphs (a permutation cycle), so we first snapshot every paragraph's
# current Range.Text into variables, then write all the new values back.
# "`v" (vertical tab, chr 11) is how Word represents a <w:br/> inside
# Range.Text, so we rebuild multi-run paragraphs using it as separator.

$d = $word.ActiveDocument

function Trim-Para([string]$s) {
    return $s.TrimEnd([char]13)
}

# --- snapshot the text currently sitting in every paragraph we touch ---
$objetivoPt  = Trim-Para $d.Paragraphs(6).Range.Text    # "Oferecer uma base sólida..."
$objetivoEn  = Trim-Para $d.Paragraphs(7).Range.Text    # "To provide a solid foundation..."
$docenteNome = Trim-Para $d.Paragraphs(9).Range.Text    # "6270264 - Juan Fernando Zapata Zapata"

# The full bilingual "Programa" content is split across two runs today
# (Portuguese at 14, English/italic at 15) built from three <w:br/>-joined
# sentences; we only need the Portuguese one, rebuilt with real separators.
$programaPt = "Números Reais: Números Naturais, Números Inteiros, Números Racionais e Irracionais. Operações com números reais, desigualdades. " + "`v" + "Funções Reais: Definição de funções, função polinomial, função racional, fatoração de polinômios, função exponencial, função logarítmica, valor absoluto, funções trigonométricas, identidades trigonométricas, funções trigonométricas inversas, funções hiperbólicas. " + "`v" + "Modelagem: Áreas, volume, custo, modelos populacionais."

$metodoVal   = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$criterioVal = "NF ≥ 5,0"
$normaVal    = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

# --- now write everything back into its new home ---

# Objetivos (Portuguese) paragraph gets the old "Programa resumido" (PT) text
$d.Paragraphs(6).Range.Text = "Números Reais, Funções Reais, modelagem com funções elementares e análise gráfica."
# Objetivos (English/italic) paragraph gets the old "Programa resumido" (EN) text
$d.Paragraphs(7).Range.Text = "Real Numbers, Real Functions, Modeling with Elementary Functions, and Graphical Analysis."

# New "Docente(s)" bullet paragraph gets the old Objetivos (PT) text
$d.Paragraphs(9).Range.Text = $objetivoPt

# "Programa resumido" (PT) paragraph gets the old full "Programa" (PT) text
$d.Paragraphs(11).Range.Text = $programaPt

# "Programa resumido" (EN/italic) paragraph gets the old Objetivos (EN) text
$d.Paragraphs(12).Range.Text = $objetivoEn

# "Programa" (PT) paragraph gets the "Método" evaluation text
$d.Paragraphs(14).Range.Text = $metodoVal

# Within the Avaliação bullet paragraph (17) there are three bold labels
# ("Método: ", "Critério: ", "Norma de recuperação: ") each followed by a
# value run; several new values collide textually with other values in the
# SAME paragraph (e.g. "NF ≥ 5,0" is both the old Critério value and the
# new Método value), so a plain Find scoped to the whole paragraph would
# hit the wrong (first) occurrence once earlier replacements land. Instead,
# scope each Find to the range AFTER the relevant bold label only.
$bibliografiaRun = "Leithold, Louis.O Cálculo com geometria Analítica: Harbra Ltda, 2009. v.1." + "`v`v" + "ANTON, Howard. Cálculo: um novo horizonte. Porto Alegre: Bookman, 2007." + "`v`v" + "THOMAS, George B. Cálculo São Paulo: Pearson Addison  Wesley, 2009. v.1," + "`v`v" + "FLEMMING, Diva M.; GONÇALVES, Mirian B. Cálculo A. São Paulo: Pearson Prentice Hall, 2009."

$p17 = $d.Paragraphs(17)

# Método: -> "NF ≥ 5,0" (paragraph-scoped Find is safe here, still only one match)
$rng = $p17.Range
$rng.Find.Execute($metodoVal, $true, $false, $false, $false, $false, $true, 1, $false, $criterioVal, 2) | Out-Null

# Critério: -> "(NF+RC)/2 ≥ 5,0, ..." — scope to text after the "Critério: " label
$labelRng = $p17.Range.Duplicate()
$labelRng.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$valRng = $d.Range($labelRng.End, $p17.Range.End)
$valRng.Find.Execute($criterioVal, $true, $false, $false, $false, $false, $true, 1, $false, $normaVal, 2) | Out-Null

# Norma de recuperação: -> bibliography text — scope to text after that label
$labelRng2 = $p17.Range.Duplicate()
$labelRng2.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$valRng2 = $d.Range($labelRng2.End, $p17.Range.End)
$valRng2.Find.Execute($normaVal, $true, $false, $false, $false, $false, $true, 1, $false, $bibliografiaRun, 2) | Out-Null

# Bibliografia paragraph (19) gets the old "Docente(s)" name bullet text
$d.Paragraphs(19).Range.Text = $docenteNome

Write-Output "done"
